$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2/E2 with new assay values, and re-round D3:E11 to 9 significant decimal digits
$ws.Range("D2").Value = 0.982952696
$ws.Range("E2").Value = 0.090851038

$ws.Range("D3").Value = 0.963293301
$ws.Range("E3").Value = 0.121493205

$ws.Range("D4").Value = 1.002794501
$ws.Range("E4").Value = 0.071071382

$ws.Range("D5").Value = 0.96978765
$ws.Range("E5").Value = 0.158575035

$ws.Range("D6").Value = 0.959784038
$ws.Range("E6").Value = 0.116524553

$ws.Range("D7").Value = 0.969431029
$ws.Range("E7").Value = 0.094286017

$ws.Range("D8").Value = 0.753018528
$ws.Range("E8").Value = 0.102221714

$ws.Range("D9").Value = 0.334670207
$ws.Range("E9").Value = 0.050067358

$ws.Range("D10").Value = 0.191020602
$ws.Range("E10").Value = 0.065054821

$ws.Range("D11").Value = 0.128803685
$ws.Range("E11").Value = 0.040077135

# Update the active selection to E6 as recorded in the saved view state
$ws.Range("E6").Select()
